# Update generated event-attendance / min-price figures on the "展览" and
# "全部类型" sheets (they mirror each other in this workbook).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> hashtable of column letter -> new value
$updates = @{
    2  = @{ F = 161; G = 35 }
    3  = @{ F = 1747 }
    8  = @{ F = 12103 }
    11 = @{ F = 484 }
    14 = @{ F = 875 }
    15 = @{ F = 13525 }
    16 = @{ F = 13578 }
    24 = @{ F = 2039 }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $cols = $updates[$row]
        foreach ($col in $cols.Keys) {
            $ws.Range("$col$row").Value = $cols[$col]
        }
    }
}
